# Fix the misaligned "Answer" column (C): the header row (row 1) had picked
# up an answer-text cell in C1 that belongs further down the column. Remove
# C1 and shift the C2:C6 answers down into C3:C7 (i.e. insert a blank cell
# at C1, pushing the column down by one row). Work from the bottom up so we
# never overwrite a source cell before it has been copied.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Copy($ws.Range("C7"))
$ws.Range("C5").Copy($ws.Range("C6"))
$ws.Range("C4").Copy($ws.Range("C5"))
$ws.Range("C3").Copy($ws.Range("C4"))
$ws.Range("C2").Copy($ws.Range("C3"))
$ws.Range("C1").Copy($ws.Range("C2"))

# C1 becomes a brand-new blank cell (no value, no style).
$ws.Range("C1").Clear()

# Restore the view: scrolled so column B is left-most, with B28:C28 selected.
$ws.Range("B28:C28").Select()
